# Fruta / hortaliza, semanal
# The underlying data rows (2-44) got reshuffled: the block of columns
# D (Fecha), I (Calidad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# P (Precio $/Kg) and Q (Kg o Unidades) for each row were moved to a new
# row position (columns A, B, C, E, F, G, H, O, R stay fixed per row).
#
# rowMap[target row] = source row (as found in the original workbook)
$rowMap = @{
    2  = 29
    3  = 15
    4  = 26
    5  = 40
    6  = 36
    7  = 25
    8  = 44
    9  = 20
    10 = 17
    11 = 28
    12 = 41
    13 = 18
    14 = 3
    15 = 22
    16 = 9
    17 = 30
    18 = 21
    19 = 6
    20 = 35
    21 = 37
    22 = 24
    23 = 19
    24 = 32
    25 = 11
    26 = 34
    27 = 16
    28 = 31
    29 = 8
    30 = 5
    31 = 23
    32 = 42
    33 = 10
    34 = 7
    35 = 12
    36 = 2
    37 = 4
    38 = 27
    39 = 14
    40 = 39
    41 = 38
    42 = 33
    43 = 43
    44 = 13
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> index
# D=4 I=9 J=10 K=11 L=12 M=13 N=14 P=16 Q=17
$cols = @(4, 9, 10, 11, 12, 13, 14, 16, 17)

# Snapshot all current values for rows 2..44 before we start overwriting,
# since the permutation has cycles longer than 2 and an in-place write
# would clobber data that is still needed.
$snapshot = @{}
for ($r = 2; $r -le 44; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Now write the values from the mapped source row into each target row.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $sourceData[$c]
    }
}
